$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; all existing data (rows 1-20) shifts down to rows 2-21
$ws.Rows.Item(1).Insert()

# New header labels
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "City Name"
$ws.Range("C1").Value = "Overnight International Visitor Spend (US$ bn)"
$ws.Range("D1").Value = "Year"

# New "Year" column for every data row
$ws.Range("D2:D21").Value = 2012

# Column C (spend) loses its currency formatting -> plain/General numbers,
# applied uniformly across the header + all data cells
$ws.Range("C1:C21").NumberFormat = "General"

# Column C is wide enough to fit its new header text
$ws.Columns.Item(3).ColumnWidth = 35.83

# Restore the selection to the newly added Year column, matching the saved view state
[void]$ws.Range("D2:D21").Select()
